$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "'2024.02.25"
$ws1.Range("C2").Value = "太仓·龙吟动漫游戏展"
$ws1.Range("D2").Value = "滨河路128号 凯景世纪大酒店(太仓滨河路店)"
$ws1.Range("E2").Value = "2024.02.25 10:00-02.25 17:00"
$ws1.Range("F2").Value = 78
$ws1.Range("G2").Value = 45
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81242"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202402/IwXBoz7t1708330463199.jpeg"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "'2024.02.25"
$ws1.Range("C3").Value = "苏州·第五届次元鹿角动漫游戏展（取消）"
$ws1.Range("D3").Value = "清禾路886号 尹山湖大剧院"
$ws1.Range("E3").Value = "2024.02.25 10:00-02.25 17:00"
$ws1.Range("F3").Value = 2707
$ws1.Range("G3").Value = "不可售"
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=79333"
$ws1.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "'2024.03.08"
$ws1.Range("C4").Value = "苏州·国风宠物-cosplay展（取消）"
$ws1.Range("D4").Value = "金山南路影视城 木渎影视城会展中心"
$ws1.Range("E4").Value = "2024.03.08 09:00-03.10 17:30"
$ws1.Range("F4").Value = 1166
$ws1.Range("G4").Value = "不可售"
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=80635"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "'2024.03.17"
$ws1.Range("C5").Value = "苏州·世纪幻想动漫游戏展2.0"
$ws1.Range("D5").Value = "清禾路886号 尹山湖大剧院"
$ws1.Range("E5").Value = "2024.03.17 10:00-03.17 17:00"
$ws1.Range("F5").Value = 1322
$ws1.Range("G5").Value = 60
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81387"
$ws1.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202402/isVyI9hH1708590817616.jpeg"

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "'2024.03.23"
$ws1.Range("C6").Value = "苏州·Look Look动漫嘉年华"
$ws1.Range("D6").Value = "阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店"
$ws1.Range("E6").Value = "2024.03.23 10:00-03.23 17:30"
$ws1.Range("F6").Value = 288
$ws1.Range("G6").Value = 58
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81698"
$ws1.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202402/ZYkvUFn41706869061984.jpeg"

$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "'2024.03.30"
$ws1.Range("C7").Value = "苏州·奇幻世界5.3动漫游戏展"
$ws1.Range("D7").Value = "龙河路1288号 乐动力苏州湾体育中心"
$ws1.Range("E7").Value = "2024.03.30 10:00-03.31 17:00"
$ws1.Range("F7").Value = 1016
$ws1.Range("G7").Value = 55
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=82002"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg"

$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "'2024.04.04"
$ws1.Range("C8").Value = "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会"
$ws1.Range("D8").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Range("E8").Value = "2024.04.04 10:00-04.05 17:00"
$ws1.Range("F8").Value = 10452
$ws1.Range("G8").Value = 60
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81827"
$ws1.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"

$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "'2024.04.05"
$ws1.Range("C9").Value = "苏州·X-party 国漫游戏嘉年华03"
$ws1.Range("D9").Value = "秋枫街与开平路交叉口西南角 爱琴海购物中心"
$ws1.Range("E9").Value = "2024.04.05 10:00-04.06 17:00"
$ws1.Range("F9").Value = 10
$ws1.Range("G9").Value = 48
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=82042"
$ws1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202402/WaQk4nUt1708679999084.jpeg"

$ws1.Range("A10").Value = 9
$ws1.Range("B10").Value = "'2024.04.06"
$ws1.Range("C10").Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
$ws1.Range("D10").Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
$ws1.Range("E10").Value = "2024.04.06 10:00-04.06 16:00"
$ws1.Range("F10").Value = 79
$ws1.Range("G10").Value = 49
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=80528"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"

$ws1.Range("A11").Value = 10
$ws1.Range("B11").Value = "'2024.04.13"
$ws1.Range("C11").Value = "苏州·绘时国乙1.0-秩序之外"
$ws1.Range("D11").Value = "石路步行街永福桥浜15号 银河广场"
$ws1.Range("E11").Value = "2024.04.13 13:30-04.13 20:00"
$ws1.Range("F11").Value = 277
$ws1.Range("G11").Value = 78
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=80789"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"

$ws1.Range("A12").Value = 11
$ws1.Range("B12").Value = "'2024.04.20"
$ws1.Range("C12").Value = "苏州·首届Redamancy动漫游戏嘉年华"
$ws1.Range("D12").Value = "清禾路886号 尹山湖大剧院"
$ws1.Range("E12").Value = "2024.04.20 10:00-04.20 17:00"
$ws1.Range("F12").Value = 1027
$ws1.Range("G12").Value = 60
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81879"
$ws1.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg"

$ws1.Range("A13").Value = 12
$ws1.Range("B13").Value = "'2024.04.21"
$ws1.Range("C13").Value = "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0"
$ws1.Range("D13").Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
$ws1.Range("E13").Value = "2024.04.21 10:00-04.21 21:00"
$ws1.Range("F13").Value = 669
$ws1.Range("G13").Value = 59.9
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=78666"
$ws1.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"

$ws1.Range("A14").Value = 13
$ws1.Range("B14").Value = "'2024.05.01"
$ws1.Range("C14").Value = "昆山·第十二届理想乡动漫游戏展"
$ws1.Range("D14").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E14").Value = "2024.05.01 10:00-05.03 17:00"
$ws1.Range("F14").Value = 11946
$ws1.Range("G14").Value = 59
$ws1.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=77196"
$ws1.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"

$ws1.Range("A15").Value = 14
$ws1.Range("B15").Value = "'2024.05.01"
$ws1.Range("C15").Value = "苏州·第十七届 I COME ACG  动漫品牌博览会"
$ws1.Range("D15").Value = "金山南路288号 广电国际会展中心"
$ws1.Range("E15").Value = "2024.05.01 10:00-05.02 17:00"
$ws1.Range("F15").Value = 12348
$ws1.Range("G15").Value = 65
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=79789"
$ws1.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"

$ws1.Range("A16").Value = 15
$ws1.Range("B16").Value = "'2024.05.02"
$ws1.Range("C16").Value = "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会"
$ws1.Range("D16").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E16").Value = "2024.05.02 14:00-05.02 16:00"
$ws1.Range("F16").Value = 30
$ws1.Range("G16").Value = 1
$ws1.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=81116"
$ws1.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"

$ws1.Range("A17").Value = 16
$ws1.Range("B17").Value = "'2024.05.02"
$ws1.Range("C17").Value = "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会"
$ws1.Range("D17").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E17").Value = "2024.05.02 14:00-05.02 16:00"
$ws1.Range("F17").Value = 111
$ws1.Range("G17").Value = 1
$ws1.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=81100"
$ws1.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"

$ws1.Range("A18").Value = 17
$ws1.Range("B18").Value = "'2024.05.02"
$ws1.Range("C18").Value = "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会"
$ws1.Range("D18").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E18").Value = "2024.05.02 14:00-05.02 16:00"
$ws1.Range("F18").Value = 17
$ws1.Range("G18").Value = 1
$ws1.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=81119"
$ws1.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"

$ws1.Range("A19").Value = 18
$ws1.Range("B19").Value = "'2024.05.03"
$ws1.Range("C19").Value = "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会"
$ws1.Range("D19").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E19").Value = "2024.05.03 14:00-05.03 16:00"
$ws1.Range("F19").Value = 27
$ws1.Range("G19").Value = 1
$ws1.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81118"
$ws1.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"

$ws1.Range("A20").Value = 19
$ws1.Range("B20").Value = "'2024.05.03"
$ws1.Range("C20").Value = "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会"
$ws1.Range("D20").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E20").Value = "2024.05.03 14:00-05.03 16:00"
$ws1.Range("F20").Value = 72
$ws1.Range("G20").Value = 1
$ws1.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81120"
$ws1.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

$ws1.Range("A21").Value = 20
$ws1.Range("B21").Value = "'2024.05.03"
$ws1.Range("C21").Value = "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会"
$ws1.Range("D21").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E21").Value = "2024.05.03 14:00-05.03 16:00"
$ws1.Range("F21").Value = 38
$ws1.Range("G21").Value = 1
$ws1.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=81114"
$ws1.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

$ws1.Range("A22:I22").Delete() | Out-Null

$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "'2024.02.25"
$ws4.Range("C2").Value = "太仓·龙吟动漫游戏展"
$ws4.Range("D2").Value = "滨河路128号 凯景世纪大酒店(太仓滨河路店)"
$ws4.Range("E2").Value = "2024.02.25 10:00-02.25 17:00"
$ws4.Range("F2").Value = 78
$ws4.Range("G2").Value = 45
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81242"
$ws4.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202402/IwXBoz7t1708330463199.jpeg"

$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "'2024.02.25"
$ws4.Range("C3").Value = "苏州·第五届次元鹿角动漫游戏展（取消）"
$ws4.Range("D3").Value = "清禾路886号 尹山湖大剧院"
$ws4.Range("E3").Value = "2024.02.25 10:00-02.25 17:00"
$ws4.Range("F3").Value = 2707
$ws4.Range("G3").Value = "不可售"
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=79333"
$ws4.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"

$ws4.Range("A4").Value = 3
$ws4.Range("B4").Value = "'2024.03.03"
$ws4.Range("C4").Value = "苏州·龙猫和他的朋友·动漫作品音乐会"
$ws4.Range("D4").Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$ws4.Range("E4").Value = "2024.03.03 19:30-03.03 21:00"
$ws4.Range("F4").Value = 9
$ws4.Range("G4").Value = 60
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81799"
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/gqnOEjvJ1707214629948.jpeg"

$ws4.Range("A5").Value = 4
$ws4.Range("B5").Value = "'2024.03.08"
$ws4.Range("C5").Value = "苏州·国风宠物-cosplay展（取消）"
$ws4.Range("D5").Value = "金山南路影视城 木渎影视城会展中心"
$ws4.Range("E5").Value = "2024.03.08 09:00-03.10 17:30"
$ws4.Range("F5").Value = 1166
$ws4.Range("G5").Value = "不可售"
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=80635"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"

$ws4.Range("A6").Value = 5
$ws4.Range("B6").Value = "'2024.03.17"
$ws4.Range("C6").Value = "苏州·世纪幻想动漫游戏展2.0"
$ws4.Range("D6").Value = "清禾路886号 尹山湖大剧院"
$ws4.Range("E6").Value = "2024.03.17 10:00-03.17 17:00"
$ws4.Range("F6").Value = 1322
$ws4.Range("G6").Value = 60
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81387"
$ws4.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202402/isVyI9hH1708590817616.jpeg"

$ws4.Range("A7").Value = 6
$ws4.Range("B7").Value = "'2024.03.23"
$ws4.Range("C7").Value = "苏州·Look Look动漫嘉年华"
$ws4.Range("D7").Value = "阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店"
$ws4.Range("E7").Value = "2024.03.23 10:00-03.23 17:30"
$ws4.Range("F7").Value = 288
$ws4.Range("G7").Value = 58
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81698"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202402/ZYkvUFn41706869061984.jpeg"

$ws4.Range("A8").Value = 7
$ws4.Range("B8").Value = "'2024.03.30"
$ws4.Range("C8").Value = "苏州·奇幻世界5.3动漫游戏展"
$ws4.Range("D8").Value = "龙河路1288号 乐动力苏州湾体育中心"
$ws4.Range("E8").Value = "2024.03.30 10:00-03.31 17:00"
$ws4.Range("F8").Value = 1016
$ws4.Range("G8").Value = 55
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82002"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg"

$ws4.Range("A9").Value = 8
$ws4.Range("B9").Value = "'2024.04.04"
$ws4.Range("C9").Value = "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会"
$ws4.Range("D9").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Range("E9").Value = "2024.04.04 10:00-04.05 17:00"
$ws4.Range("F9").Value = 10452
$ws4.Range("G9").Value = 60
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81827"
$ws4.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"

$ws4.Range("A10").Value = 9
$ws4.Range("B10").Value = "'2024.04.05"
$ws4.Range("C10").Value = "苏州·X-party 国漫游戏嘉年华03"
$ws4.Range("D10").Value = "秋枫街与开平路交叉口西南角 爱琴海购物中心"
$ws4.Range("E10").Value = "2024.04.05 10:00-04.06 17:00"
$ws4.Range("F10").Value = 10
$ws4.Range("G10").Value = 48
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82042"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202402/WaQk4nUt1708679999084.jpeg"

$ws4.Range("A11").Value = 10
$ws4.Range("B11").Value = "'2024.04.06"
$ws4.Range("C11").Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
$ws4.Range("D11").Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
$ws4.Range("E11").Value = "2024.04.06 10:00-04.06 16:00"
$ws4.Range("F11").Value = 79
$ws4.Range("G11").Value = 49
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=80528"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"

$ws4.Range("A12").Value = 11
$ws4.Range("B12").Value = "'2024.04.13"
$ws4.Range("C12").Value = "苏州·绘时国乙1.0-秩序之外"
$ws4.Range("D12").Value = "石路步行街永福桥浜15号 银河广场"
$ws4.Range("E12").Value = "2024.04.13 13:30-04.13 20:00"
$ws4.Range("F12").Value = 277
$ws4.Range("G12").Value = 78
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=80789"
$ws4.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"

$ws4.Range("A13").Value = 12
$ws4.Range("B13").Value = "'2024.04.20"
$ws4.Range("C13").Value = "苏州·首届Redamancy动漫游戏嘉年华"
$ws4.Range("D13").Value = "清禾路886号 尹山湖大剧院"
$ws4.Range("E13").Value = "2024.04.20 10:00-04.20 17:00"
$ws4.Range("F13").Value = 1027
$ws4.Range("G13").Value = 60
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81879"
$ws4.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg"

$ws4.Range("A14").Value = 13
$ws4.Range("B14").Value = "'2024.04.21"
$ws4.Range("C14").Value = "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0"
$ws4.Range("D14").Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
$ws4.Range("E14").Value = "2024.04.21 10:00-04.21 21:00"
$ws4.Range("F14").Value = 669
$ws4.Range("G14").Value = 59.9
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=78666"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"

$ws4.Range("A15").Value = 14
$ws4.Range("B15").Value = "'2024.05.01"
$ws4.Range("C15").Value = "昆山·第十二届理想乡动漫游戏展"
$ws4.Range("D15").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E15").Value = "2024.05.01 10:00-05.03 17:00"
$ws4.Range("F15").Value = 11946
$ws4.Range("G15").Value = 59
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=77196"
$ws4.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"

$ws4.Range("A16").Value = 15
$ws4.Range("B16").Value = "'2024.05.01"
$ws4.Range("C16").Value = "苏州·第十七届 I COME ACG  动漫品牌博览会"
$ws4.Range("D16").Value = "金山南路288号 广电国际会展中心"
$ws4.Range("E16").Value = "2024.05.01 10:00-05.02 17:00"
$ws4.Range("F16").Value = 12348
$ws4.Range("G16").Value = 65
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=79789"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"

$ws4.Range("A17").Value = 16
$ws4.Range("B17").Value = "'2024.05.02"
$ws4.Range("C17").Value = "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会"
$ws4.Range("D17").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E17").Value = "2024.05.02 14:00-05.02 16:00"
$ws4.Range("F17").Value = 30
$ws4.Range("G17").Value = 1
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=81116"
$ws4.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"

$ws4.Range("A18").Value = 17
$ws4.Range("B18").Value = "'2024.05.02"
$ws4.Range("C18").Value = "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会"
$ws4.Range("D18").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E18").Value = "2024.05.02 14:00-05.02 16:00"
$ws4.Range("F18").Value = 111
$ws4.Range("G18").Value = 1
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=81100"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"

$ws4.Range("A19").Value = 18
$ws4.Range("B19").Value = "'2024.05.02"
$ws4.Range("C19").Value = "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会"
$ws4.Range("D19").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E19").Value = "2024.05.02 14:00-05.02 16:00"
$ws4.Range("F19").Value = 17
$ws4.Range("G19").Value = 1
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81119"
$ws4.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"

$ws4.Range("A20").Value = 19
$ws4.Range("B20").Value = "'2024.05.03"
$ws4.Range("C20").Value = "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会"
$ws4.Range("D20").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E20").Value = "2024.05.03 14:00-05.03 16:00"
$ws4.Range("F20").Value = 27
$ws4.Range("G20").Value = 1
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81118"
$ws4.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"

$ws4.Range("A21").Value = 20
$ws4.Range("B21").Value = "'2024.05.03"
$ws4.Range("C21").Value = "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会"
$ws4.Range("D21").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E21").Value = "2024.05.03 14:00-05.03 16:00"
$ws4.Range("F21").Value = 72
$ws4.Range("G21").Value = 1
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=81120"
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"

$ws4.Range("A22").Value = 21
$ws4.Range("B22").Value = "'2024.05.03"
$ws4.Range("C22").Value = "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会"
$ws4.Range("D22").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E22").Value = "2024.05.03 14:00-05.03 16:00"
$ws4.Range("F22").Value = 38
$ws4.Range("G22").Value = 1
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=81114"
$ws4.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

$ws4.Range("A23:I23").Delete() | Out-Null

